$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 62.58
$ws.Range("I15").Value = 62.58
$ws.Range("K15").Value = 187.74
$ws.Range("M15").Value = -18.74000000000001
$ws.Range("H17").Value = 3050566.2
$ws.Range("J17").Value = 3145837
$ws.Range("L17").Value = 9437511
$ws.Range("N17").Value = -9437847
$ws.Range("H76").Value = 3706559.8
$ws.Range("I76").Value = 3077
$ws.Range("J76").Value = 6947107
$ws.Range("K76").Value = 3077
$ws.Range("L76").Value = 6947107
$ws.Range("M76").Value = -2762
$ws.Range("N76").Value = -6947737
$ws.Range("H79").Value = 3706559.8
$ws.Range("I79").Value = 3077
$ws.Range("J79").Value = 6947107
$ws.Range("K79").Value = 3077
$ws.Range("L79").Value = 6947107
$ws.Range("M79").Value = -1985
$ws.Range("N79").Value = -6949291
$ws.Range("H111").Value = 3150.7144
$ws.Range("I111").Value = 4019.75
$ws.Range("J111").Value = 1992
$ws.Range("K111").Value = 12059.25
$ws.Range("L111").Value = 5976
$ws.Range("M111").Value = -8992.25
$ws.Range("N111").Value = -12110
$ws.Range("H132").Value = 32261202
$ws.Range("I132").Value = 35717636
$ws.Range("K132").Value = 107152908
$ws.Range("M132").Value = -107150378

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8967.857
$ws.Range("I32").Value = 6597.5137
$ws.Range("J32").Value = 23189.916
$ws.Range("K32").Value = 6597.5137
$ws.Range("L32").Value = 23189.916
$ws.Range("M32").Value = -6310.5137
$ws.Range("N32").Value = -23763.916
$ws.Range("H61").Value = 10755067
$ws.Range("I61").Value = 13335107
$ws.Range("J61").Value = 4899.3335
$ws.Range("K61").Value = 13335107
$ws.Range("L61").Value = 4899.3335
$ws.Range("M61").Value = -13334895
$ws.Range("N61").Value = -5323.3335
$ws.Range("H74").Value = 47620800
$ws.Range("I74").Value = 62500760
$ws.Range("J74").Value = 4920
$ws.Range("K74").Value = 62500760
$ws.Range("L74").Value = 4920
$ws.Range("M74").Value = -62499886
$ws.Range("N74").Value = -6668
$ws.Range("H77").Value = 47620800
$ws.Range("I77").Value = 62500760
$ws.Range("J77").Value = 4920
$ws.Range("K77").Value = 312503800
$ws.Range("L77").Value = 24600
$ws.Range("M77").Value = -312499432
$ws.Range("N77").Value = -33336
$ws.Range("H132").Value = 9815890
$ws.Range("I132").Value = 11906660
$ws.Range("K132").Value = 35719980
$ws.Range("M132").Value = -35717450
$ws.Range("H136").Value = 10755067
$ws.Range("I136").Value = 13335107
$ws.Range("J136").Value = 4899.3335
$ws.Range("K136").Value = 40005321
$ws.Range("L136").Value = 14698.0005
$ws.Range("M136").Value = -40002771
$ws.Range("N136").Value = -19798.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3479.2
$ws.Range("I134").Value = 3479.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10437.6
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7902.599999999999
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 50524
$ws.Range("I36").Value = 50524
$ws.Range("K36").Value = 50524
$ws.Range("M36").Value = -50136
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H40").Value = 50524
$ws.Range("I40").Value = 50524
$ws.Range("K40").Value = 50524
$ws.Range("M40").Value = -50364
$ws.Range("H42").Value = 6000
$ws.Range("I42").Value = 6000
$ws.Range("K42").Value = 6000
$ws.Range("M42").Value = -5407
$ws.Range("H58").Value = 15706.114
$ws.Range("I58").Value = 1440.579
$ws.Range("J58").Value = 32646.438
$ws.Range("K58").Value = 1440.579
$ws.Range("L58").Value = 32646.438
$ws.Range("M58").Value = -1237.579
$ws.Range("N58").Value = -33052.43799999999
$ws.Range("H132").Value = 50003130
$ws.Range("I132").Value = 58825612
$ws.Range("K132").Value = 176476836
$ws.Range("M132").Value = -176474306
$ws.Range("H136").Value = 15706.114
$ws.Range("I136").Value = 1440.579
$ws.Range("J136").Value = 32646.438
$ws.Range("K136").Value = 4321.737
$ws.Range("L136").Value = 97939.314
$ws.Range("M136").Value = -1771.737
$ws.Range("N136").Value = -103039.314

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1286.7174
$ws.Range("I5").Value = 882.38464
$ws.Range("J5").Value = 1812.35
$ws.Range("K5").Value = 2647.15392
$ws.Range("L5").Value = 5437.049999999999
$ws.Range("M5").Value = -2535.15392
$ws.Range("N5").Value = -5661.049999999999
$ws.Range("H112").Value = 125001480
$ws.Range("J112").Value = 500002500
$ws.Range("L112").Value = 1500007500
$ws.Range("N112").Value = -1500009716
$ws.Range("H115").Value = 4976.3
$ws.Range("I115").Value = 1256.6666
$ws.Range("J115").Value = 6570.4287
$ws.Range("K115").Value = 3769.9998
$ws.Range("L115").Value = 19711.2861
$ws.Range("M115").Value = -2594.9998
$ws.Range("N115").Value = -22061.2861
$ws.Range("H118").Value = 250000130
$ws.Range("I118").Value = 250000130
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 750000390
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -749999147
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 1055.4166
$ws.Range("J121").Value = 1127.2727
$ws.Range("L121").Value = 3381.8181
$ws.Range("N121").Value = -6001.8181
$ws.Range("H122").Value = 1296
$ws.Range("J122").Value = 1296
$ws.Range("L122").Value = 11664
$ws.Range("N122").Value = -16564
$ws.Range("H131").Value = 687.78
$ws.Range("J131").Value = 733.6667
$ws.Range("L131").Value = 2201.0001
$ws.Range("N131").Value = -12281.0001
$ws.Range("H134").Value = 2806.457
$ws.Range("I134").Value = 1858.25
$ws.Range("J134").Value = 6599.2856
$ws.Range("K134").Value = 5574.75
$ws.Range("L134").Value = 19797.8568
$ws.Range("M134").Value = -504.75
$ws.Range("N134").Value = -29937.8568
$ws.Range("H135").Value = 1286.7174
$ws.Range("I135").Value = 882.38464
$ws.Range("J135").Value = 1812.35
$ws.Range("K135").Value = 7941.46176
$ws.Range("L135").Value = 16311.15
$ws.Range("M135").Value = -5406.46176
$ws.Range("N135").Value = -21381.15

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3754056.2
$ws.Range("I132").Value = 5525822.5
$ws.Range("J132").Value = 49454.273
$ws.Range("K132").Value = 16577467.5
$ws.Range("L132").Value = 148362.819
$ws.Range("M132").Value = -16574937.5
$ws.Range("N132").Value = -153422.819

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2896.9443
$ws.Range("J7").Value = 2802.7144
$ws.Range("L7").Value = 2802.7144
$ws.Range("N7").Value = -3026.7144
$ws.Range("H126").Value = 2896.9443
$ws.Range("J126").Value = 2802.7144
$ws.Range("L126").Value = 8408.143199999999
$ws.Range("N126").Value = -13348.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10000400
$ws.Range("J11").Value = 800
$ws.Range("L11").Value = 800
$ws.Range("N11").Value = -1084
$ws.Range("H132").Value = 15152731
$ws.Range("I132").Value = 20000736
$ws.Range("J132").Value = 2719.125
$ws.Range("K132").Value = 60002208
$ws.Range("L132").Value = 8157.375
$ws.Range("M132").Value = -59999678
$ws.Range("N132").Value = -13217.375
